# Applies the cryptos.xlsx price/volume refresh described in the diff.
# For numeric-looking Price (column D) values we prefix with a leading
# apostrophe so Excel stores them as literal text (matching the source
# workbook, which keeps every Price/Volume cell as a text string) instead
# of silently parsing them into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "92.727.06"
$ws.Range("E2").Value = "  -1.99%  "

$ws.Range("D3").Value = "3.355.14"
$ws.Range("E3").Value = "  -2.55%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'232.05"
$ws.Range("E5").Value = "  -2.35%  "

$ws.Range("D6").Value = "'616.57"
$ws.Range("E6").Value = "  -4.19%  "

$ws.Range("D7").Value = "'1.36"
$ws.Range("E7").Value = "  -5.58%  "

$ws.Range("D8").Value = "'0.387"
$ws.Range("E8").Value = "  -4.76%  "

$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Value = "'0.944"
$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("D11").Value = "3.354.95"
$ws.Range("E11").Value = "  -2.48%  "

$ws.Range("D12").Value = "'42.51"
$ws.Range("E12").Value = "  +1.72%  "

$ws.Range("D13").Value = "'0.196"
$ws.Range("E13").Value = "  -1.16%  "

$ws.Range("D14").Value = "'6.21"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").Value = "92.678.15"
$ws.Range("E15").Value = "  -1.73%  "

$ws.Range("D16").Value = "3.991.22"
$ws.Range("E16").Value = "  -2.08%  "

$ws.Range("D17").Value = "'0.0000244"
$ws.Range("E17").Value = "  -3.05%  "

$ws.Range("D18").Value = "'8.04"
$ws.Range("E18").Value = "  -3.68%  "

$ws.Range("D19").Value = "3.346.22"
$ws.Range("E19").Value = "  -2.73%  "

$ws.Range("D20").Value = "'17.31"
$ws.Range("E20").Value = "  -1.53%  "

$ws.Range("D21").Value = "'11.20"
$ws.Range("E21").Value = "  -2.72%  "

$ws.Range("D22").Value = "'3.33"
$ws.Range("E22").Value = "  +3.29%  "

$ws.Range("D23").Value = "'493.40"
$ws.Range("E23").Value = "  -1.48%  "

$ws.Range("D24").Value = "'0.427"
$ws.Range("E24").Value = "  -15.70%  "

$ws.Range("D25").Value = "'6.60"
$ws.Range("E25").Value = "  +1.31%  "

$ws.Range("D26").Value = "'0.0000182"
$ws.Range("E26").Value = "  -5.65%  "

$ws.Range("D27").Value = "'92.68"
$ws.Range("E27").Value = "  -1.61%  "

$ws.Range("D28").Value = "'11.94"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").Value = "3.535.57"
$ws.Range("E29").Value = "  -2.33%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").Value = "'11.02"
$ws.Range("E31").Value = "  -6.15%  "

$ws.Range("D32").Value = "'0.135"
$ws.Range("E32").Value = "  -2.57%  "

$ws.Range("D33").Value = "'2.66"
$ws.Range("E33").Value = "  -3.51%  "

$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("D35").Value = "'0.171"
$ws.Range("E35").Value = "  -4.30%  "

$ws.Range("D36").Value = "'28.44"
$ws.Range("E36").Value = "  -5.10%  "

$ws.Range("D37").Value = "'0.524"
$ws.Range("E37").Value = "  -5.56%  "

$ws.Range("D38").Value = "'551.51"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").Value = "'7.43"
$ws.Range("E39").Value = "  -3.13%  "

$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("E41").Value = "  -1.74%  "

$ws.Range("E42").Value = "  -5.00%  "

$ws.Range("D43").Value = "'0.877"
$ws.Range("E43").Value = "  -3.71%  "

$ws.Range("D44").Value = "'23.65"
$ws.Range("E44").Value = "  -1.68%  "

$ws.Range("D45").Value = "'1.69"
$ws.Range("E45").Value = "  -1.02%  "

$ws.Range("D46").Value = "'3.57"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("D47").Value = "'0.0405"
$ws.Range("E47").Value = "  -1.38%  "

$ws.Range("D48").Value = "'5.37"
$ws.Range("E48").Value = "  -3.93%  "

$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").Value = "'52.57"
$ws.Range("E49").Value = "  -3.45%  "

$ws.Range("B50").Value = "Fantom"
$ws.Range("C50").Value = "https://coinranking.com/coin/uIEWfMFnQo9K_+fantom-ftm"
$ws.Range("D50").Value = "'1.12"
$ws.Range("E50").Value = "  +17.09%  "

$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'2.09"
$ws.Range("E51").Value = "  -3.81%  "
